# "Generate Report for Handoff"
#
# For the rows that are "Ready for handoff" (rows 4-7 on both the zh-cn and
# de-de localization-status sheets), refresh the report:
#   - Priority moves from "low" to "ht" (a new handoff pass was generated)
#   - Latest Handoff Datetime is bumped to the new generation timestamp
#
# zh-cn rows 4-7: 2016-08-21 14:40:10 -> 2016-08-21 14:40:26
# de-de rows 4-7: 2016-08-21 14:40:14 -> 2016-08-21 14:40:30

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($row in 4..7) {
    $zhcn.Cells.Item($row, 5).Value = "ht"
    $zhcn.Cells.Item($row, 8).Value = "2016-08-21 14:40:26"
}

$dede = $wb.Worksheets.Item("de-de")
foreach ($row in 4..7) {
    $dede.Cells.Item($row, 5).Value = "ht"
    $dede.Cells.Item($row, 8).Value = "2016-08-21 14:40:30"
}
